$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 1).Value = 112044333
$ws.Cells.Item(3, 2).Value = 56398
$ws.Cells.Item(3, 3).Value = 'Ovaliderad'
$ws.Cells.Item(3, 4).Value = 'NT'
$ws.Cells.Item(3, 5).Value = 100109
$ws.Cells.Item(3, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(3, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(3, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(3, 16).Value = 'Stor Mpmerg, Kilen-Stor, Moberg, Leksand, Dlr'
$ws.Cells.Item(3, 17).Value = 511613.7990622812
$ws.Cells.Item(3, 18).Value = 6733639.811082688
$ws.Cells.Item(3, 19).Value = 25
$ws.Cells.Item(3, 20).Value = 'Dalarna'
$ws.Cells.Item(3, 21).Value = 'Leksand'
$ws.Cells.Item(3, 22).Value = 'Dalarna'
$ws.Cells.Item(3, 23).Value = 'Leksand'
$ws.Cells.Item(3, 25).NumberFormat = "@"
$ws.Cells.Item(3, 25).Value = '2023-09-12'
$ws.Cells.Item(3, 25).Style = "Normal"
$ws.Cells.Item(3, 26).Value = '00:00'
$ws.Cells.Item(3, 27).NumberFormat = "@"
$ws.Cells.Item(3, 27).Value = '2023-09-12'
$ws.Cells.Item(3, 27).Style = "Normal"
$ws.Cells.Item(3, 28).Value = '00:00'
$ws.Cells.Item(3, 29).Value = 'Minst 2'
$ws.Cells.Item(3, 30).Value = $false
$ws.Cells.Item(3, 31).Value = $false
$ws.Cells.Item(3, 33).Value = $false
$ws.Cells.Item(3, 49).Value = 'Åke Sköld'
$ws.Cells.Item(3, 50).Value = 'Åke Sköld'

# Row 4
$ws.Cells.Item(4, 1).Value = 112043031
$ws.Cells.Item(4, 2).Value = 90332
$ws.Cells.Item(4, 3).Value = 'Ovaliderad'
$ws.Cells.Item(4, 4).Value = 'LC'
$ws.Cells.Item(4, 5).Value = 4769
$ws.Cells.Item(4, 6).Value = 'Svavelriska'
$ws.Cells.Item(4, 7).Value = 'Lactarius scrobiculatus'
$ws.Cells.Item(4, 8).Value = '(Scop.:Fr.) Fr.'
$ws.Cells.Item(4, 16).Value = 'Stor-Moberg (Stor-Moberg), Dlr'
$ws.Cells.Item(4, 17).Value = 511625.1419049087
$ws.Cells.Item(4, 18).Value = 6733616.372369035
$ws.Cells.Item(4, 19).Value = 1
$ws.Cells.Item(4, 20).Value = 'Dalarna'
$ws.Cells.Item(4, 21).Value = 'Leksand'
$ws.Cells.Item(4, 22).Value = 'Dalarna'
$ws.Cells.Item(4, 23).Value = 'Leksand'
$ws.Cells.Item(4, 25).NumberFormat = "@"
$ws.Cells.Item(4, 25).Value = '2023-09-12'
$ws.Cells.Item(4, 25).Style = "Normal"
$ws.Cells.Item(4, 26).Value = '10:42'
$ws.Cells.Item(4, 27).NumberFormat = "@"
$ws.Cells.Item(4, 27).Value = '2023-09-12'
$ws.Cells.Item(4, 27).Style = "Normal"
$ws.Cells.Item(4, 28).Value = '10:42'
$ws.Cells.Item(4, 30).Value = $false
$ws.Cells.Item(4, 31).Value = $false
$ws.Cells.Item(4, 33).Value = $false
$ws.Cells.Item(4, 49).Value = 'Evalena Sköld'
$ws.Cells.Item(4, 50).Value = 'Evalena Sköld'

# Row 5
$ws.Cells.Item(5, 1).Value = 112042940
$ws.Cells.Item(5, 2).Value = 98535
$ws.Cells.Item(5, 3).Value = 'Ovaliderad'
$ws.Cells.Item(5, 4).Value = 'LC'
$ws.Cells.Item(5, 5).Value = 222498
$ws.Cells.Item(5, 6).Value = 'Blåsippa'
$ws.Cells.Item(5, 7).Value = 'Hepatica nobilis'
$ws.Cells.Item(5, 8).Value = 'Schreb.'
$ws.Cells.Item(5, 16).Value = 'Stor-Moberg (Stor-Moberg), Dlr'
$ws.Cells.Item(5, 17).Value = 511610.9043343531
$ws.Cells.Item(5, 18).Value = 6733626.107665217
$ws.Cells.Item(5, 19).Value = 1
$ws.Cells.Item(5, 20).Value = 'Dalarna'
$ws.Cells.Item(5, 21).Value = 'Leksand'
$ws.Cells.Item(5, 22).Value = 'Dalarna'
$ws.Cells.Item(5, 23).Value = 'Leksand'
$ws.Cells.Item(5, 25).NumberFormat = "@"
$ws.Cells.Item(5, 25).Value = '2023-09-12'
$ws.Cells.Item(5, 25).Style = "Normal"
$ws.Cells.Item(5, 26).Value = '10:33'
$ws.Cells.Item(5, 27).NumberFormat = "@"
$ws.Cells.Item(5, 27).Value = '2023-09-12'
$ws.Cells.Item(5, 27).Style = "Normal"
$ws.Cells.Item(5, 28).Value = '10:33'
$ws.Cells.Item(5, 29).Value = 'Fullt med blåsippsblad på denna sidan bäcken'
$ws.Cells.Item(5, 30).Value = $false
$ws.Cells.Item(5, 31).Value = $false
$ws.Cells.Item(5, 33).Value = $false
$ws.Cells.Item(5, 49).Value = 'Evalena Sköld'
$ws.Cells.Item(5, 50).Value = 'Evalena Sköld, Åke Sköld'

# Row 6
$ws.Cells.Item(6, 1).Value = 112043158
$ws.Cells.Item(6, 2).Value = 95532
$ws.Cells.Item(6, 3).Value = 'Ovaliderad'
$ws.Cells.Item(6, 4).Value = 'LC'
$ws.Cells.Item(6, 5).Value = 221945
$ws.Cells.Item(6, 6).Value = 'Revlummer'
$ws.Cells.Item(6, 7).Value = 'Lycopodium annotinum'
$ws.Cells.Item(6, 8).Value = 'L.'
$ws.Cells.Item(6, 16).Value = 'Stor-Moberg (Stor-Moberg), Dlr'
$ws.Cells.Item(6, 17).Value = 511628.0588172724
$ws.Cells.Item(6, 18).Value = 6733623.228879539
$ws.Cells.Item(6, 19).Value = 1
$ws.Cells.Item(6, 20).Value = 'Dalarna'
$ws.Cells.Item(6, 21).Value = 'Leksand'
$ws.Cells.Item(6, 22).Value = 'Dalarna'
$ws.Cells.Item(6, 23).Value = 'Leksand'
$ws.Cells.Item(6, 25).NumberFormat = "@"
$ws.Cells.Item(6, 25).Value = '2023-09-12'
$ws.Cells.Item(6, 25).Style = "Normal"
$ws.Cells.Item(6, 26).Value = '10:51'
$ws.Cells.Item(6, 27).NumberFormat = "@"
$ws.Cells.Item(6, 27).Value = '2023-09-12'
$ws.Cells.Item(6, 27).Style = "Normal"
$ws.Cells.Item(6, 28).Value = '10:51'
$ws.Cells.Item(6, 29).Value = 'Finns fläckvis i området'
$ws.Cells.Item(6, 30).Value = $false
$ws.Cells.Item(6, 31).Value = $false
$ws.Cells.Item(6, 33).Value = $false
$ws.Cells.Item(6, 49).Value = 'Evalena Sköld'
$ws.Cells.Item(6, 50).Value = 'Evalena Sköld, Åke Sköld'
